{"js": "// The final paragraph ends with the run text \"\u4e09\u677f\u3002\" (part of\n// \"\u7535\u529b\u677f\u5757\u957f\u57ce\u7535\u5de5\u4e09\u677f\u3002\"). The edit expands that into:\n//   \"\u4e09\u677f\uff0c\u5d07\u5fb7\u79d1\u6280\u5f3a\u8d8b\u52bf \u7535\u529b\u8bbe\u5907\u53e0\u52a0\u6b21\u65b0\u3002\u7f8e\u829d\u80a1\u4efd\u5730\u4ea7\u94fe\u88c5\u4fee\u88c5\u9970\u3002\u65b0\u6d01\u80fd\u3001\u5bd2\u6b66\u7eaa\u8d8b\u52bf\u65b0\u9ad8\u3002\"\n// while keeping the existing run formatting (plain text, eastAsia font hint).\nconst body = context.document.body;\n\nconst results = body.search(\"\u4e09\u677f\u3002\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text '\u4e09\u677f\u3002' not found in document body.\");\n}\n\n// Use the last match in case of repeats; in this document there is exactly one.\nconst target = results.items[results.items.length - 1];\n\nconst newText =\n  \"\u4e09\u677f\" +\n  \"\uff0c\u5d07\u5fb7\u79d1\u6280\u5f3a\u8d8b\u52bf \u7535\u529b\u8bbe\u5907\u53e0\u52a0\u6b21\u65b0\" +\n  \"\u3002\" +\n  \"\u7f8e\u829d\u80a1\u4efd\u5730\u4ea7\u94fe\u88c5\u4fee\u88c5\u9970\u3002\u65b0\u6d01\u80fd\u3001\u5bd2\u6b66\u7eaa\u8d8b\u52bf\u65b0\u9ad8\u3002\";\n\ntarget.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# The last paragraph of the document ends in a run containing \"\u4e09\u677f\u3002\"\n# (part of \"\u7535\u529b\u677f\u5757\u957f\u57ce\u7535\u5de5\u4e09\u677f\u3002\"). Expand that run's text into:\n#   \"\u4e09\u677f\uff0c\u5d07\u5fb7\u79d1\u6280\u5f3a\u8d8b\u52bf \u7535\u529b\u8bbe\u5907\u53e0\u52a0\u6b21\u65b0\u3002\u7f8e\u829d\u80a1\u4efd\u5730\u4ea7\u94fe\u88c5\u4fee\u88c5\u9970\u3002\u65b0\u6d01\u80fd\u3001\u5bd2\u6b66\u7eaa\u8d8b\u52bf\u65b0\u9ad8\u3002\"\n# keeping the surrounding (unchanged) run formatting.\n\n$d = $word.ActiveDocument\n\n$newText = \"\u4e09\u677f\" + \"\uff0c\u5d07\u5fb7\u79d1\u6280\u5f3a\u8d8b\u52bf \u7535\u529b\u8bbe\u5907\u53e0\u52a0\u6b21\u65b0\" + \"\u3002\" + \"\u7f8e\u829d\u80a1\u4efd\u5730\u4ea7\u94fe\u88c5\u4fee\u88c5\u9970\u3002\u65b0\u6d01\u80fd\u3001\u5bd2\u6b66\u7eaa\u8d8b\u52bf\u65b0\u9ad8\u3002\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"\u4e09\u677f\u3002\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith,\n# Replace(wdReplaceAll=2)\n$find.Execute(\"\u4e09\u677f\u3002\", $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n"}
